$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Previously added")
$ws2 = $wb.Worksheets.Item("New")

# Remember the hyperlink target URLs on the "New" sheet (rows 2-5, col A)
# before we touch anything, since we'll need to recreate them on the
# "Previously added" sheet. In this workbook the displayed cell text IS
# the link target, so read it off the cell (Hyperlink.Address doesn't
# reliably round-trip through this host).
$links = @()
for ($i = 2; $i -le 5; $i++) {
    $links += $ws2.Cells.Item($i, 1).Value2
}

# Find the first free row at the bottom of "Previously added".
$lastRow = $ws1.UsedRange.Rows.Count
$destFirstRow = $lastRow + 1

# Pre-create the hyperlinks on the (currently empty) destination cells.
# Hyperlinks.Add stamps its own "Hyperlink" cell style, so do this BEFORE
# copying the real row data/format over them - the subsequent Copy()
# overwrites the cell format (restoring the original look) while leaving
# the already-attached hyperlink relationship in place.
for ($i = 0; $i -lt 4; $i++) {
    $destRow = $destFirstRow + $i
    $ws1.Hyperlinks.Add($ws1.Cells.Item($destRow, 1), $links[$i])
}

# Copy the 4 data rows (A2:F5) from "New" down to the bottom of
# "Previously added", preserving values + formatting.
$ws2.Range("A2:F5").Copy($ws1.Range("A" + $destFirstRow))

# Remove the now-duplicated hyperlinks from "New" and clear its data rows,
# leaving only the header row behind.
$ws2.Hyperlinks.Delete()
$ws2.Range("A2:F5").Clear()
